$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates (Target cluster stays ECs) ---
$ws.Range("G2").Value = 0.0345785
$ws.Range("H2").Value = 0.069157
$ws.Range("M2").Value = 8.961497999999999
$ws.Range("N2").Value = 26.884494
$ws.Range("O2").Value = 0.1368263930819497
$ws.Range("P2").Value = 0.1434363050389826
$ws.Range("Q2").Value = 0.309875158593
$ws.Range("R2").Value = 1.859250951558
$ws.Range("S2").Value = 0.1368263930819497
$ws.Range("T2").Value = 0.1434363050389826

# --- Row 3 updates (Target cluster stays FAPs) ---
$ws.Range("G3").Value = 0.0345785
$ws.Range("H3").Value = 0.069157
$ws.Range("O3").Value = 0.7230214505096683
$ws.Range("P3").Value = 0.7579497126911668
$ws.Range("Q3").Value = 1.637450067902
$ws.Range("R3").Value = 9.824700407411999
$ws.Range("S3").Value = 0.7230214505096683
$ws.Range("T3").Value = 0.7579497126911668

# --- Row 4 updates (Target cluster changes from MuSCs to Inflammatory-Mac) ---
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("G4").Value = 0.0345785
$ws.Range("H4").Value = 0.069157
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.110461
$ws.Range("N4").Value = 0.331383
$ws.Range("O4").Value = 0.00168654617857698
$ws.Range("P4").Value = 0.00176802111554464
$ws.Range("Q4").Value = 0.0038195756885
$ws.Range("R4").Value = 0.022917454131
$ws.Range("S4").Value = 0.00168654617857698
$ws.Range("T4").Value = 0.00176802111554464

# --- New Row 5 (Target cluster MuSCs) ---
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Plg"
$ws.Range("C5").Value = "F2r"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.5
$ws.Range("G5").Value = 0.0345785
$ws.Range("H5").Value = 0.069157
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 9.0545855
$ws.Range("N5").Value = 18.109171
$ws.Range("O5").Value = 0.1382476763167411
$ws.Range("P5").Value = 0.09661749912641458
$ws.Range("Q5").Value = 0.31309398471175
$ws.Range("R5").Value = 1.252375938847
$ws.Range("S5").Value = 0.1382476763167411
$ws.Range("T5").Value = 0.09661749912641458

# --- New Row 6 (Target cluster Resolving-Mac) ---
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Plg"
$ws.Range("C6").Value = "F2r"
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.5
$ws.Range("G6").Value = 0.0345785
$ws.Range("H6").Value = 0.069157
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 1
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.01427366666666667
$ws.Range("N6").Value = 0.042821
$ws.Range("O6").Value = 0.0002179339130638713
$ws.Range("P6").Value = 0.0002284620278914037
$ws.Range("Q6").Value = 0.0004935619828333333
$ws.Range("R6").Value = 0.002961371897
$ws.Range("S6").Value = 0.0002179339130638713
$ws.Range("T6").Value = 0.0002284620278914037
